$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Add new row of data to Sheet1 (row 27: "Pascal's Triangle") ---
# Add the hyperlink first (TextToDisplay gives the <hyperlink display="..."/> attribute,
# matching the URL text), then overwrite the cell value with the real problem title -
# this keeps the friendly cell text without losing the display attribute.
$ws1.Hyperlinks.Add($ws1.Range("B27"), "https://leetcode.com/problems/pascals-triangle/", "", "", "https://leetcode.com/problems/pascals-triangle/")
$ws1.Range("B27").Value = "Pascal's Triangle"
$ws1.Range("B27").Style = "Hyperlink"

$ws1.Range("C27").Value = 1
$ws1.Range("D27").Value = 1
$ws1.Range("E27").Value = 38
$ws1.Range("F27").Value = 0.25
$ws1.Range("G27").Value = 16.3
$ws1.Range("H27").Value = 0.69
$ws1.Range("I27").Value = "https://leetcode.com/problems/pascals-triangle/submissions/1068158329/"

# --- Update sheet view / selection state ---
$ws1.Activate() | Out-Null
$ws1.Range("I27").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 2

$ws2.Activate() | Out-Null
$ws2.Range("E6").Select() | Out-Null
